# Insert a new data row at row 43, pushing all existing rows (43-134) down
# by one (to 44-135), and populate the newly inserted row with the new
# observation described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 43 (shifts cells down).
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with the new record's data.
$ws.Range("A43").Value = 2
$ws.Range("B43").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44868
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = 100112024
$ws.Range("G43").Value = "Choclo"
$ws.Range("H43").Value = "Dulce o Americano"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 700
$ws.Range("K43").Value = 30000
$ws.Range("L43").Value = 32000
$ws.Range("M43").Value = 31000
$ws.Range("N43").Value = "$/malla 70 unidades"
$ws.Range("O43").Value = "Provincia de Limarí"
$ws.Range("P43").Value = 443
$ws.Range("Q43").Value = 70
$ws.Range("R43").Value = "Hortaliza"
